$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("D2").Value = 300
$ws.Range("D3").Value = 200
$ws.Range("D4").Value = 2000
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 400
$ws.Range("D6").Value = 1700
$ws.Range("D7").Value = 1000
$ws.Range("D8").Value = 1200
$ws.Range("D9").Value = 600
$ws.Range("D10").Value = 400
$ws.Range("D11").Value = 300

$ws.Protect()
